$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Manchester tribunal block (rows 3-7) ---
$ws.Range("A3").Value = "tribunalManchesterAddress"
$ws.Range("B3").Value = "35 La Nava S3 6AD, Southampton"

$ws.Range("A4").Value = "tribunalManchesterTelephone"
$ws.Range("B4").Value = 3577131270

$ws.Range("A5").Value = "tribunalManchesterFax"
$ws.Range("B5").Value = 7577126570

$ws.Range("A6").Value = "tribunalManchesterDX"
$ws.Range("B6").Value = 123456

$ws.Range("A7").Value = "tribunalManchesterEmail"
$ws.Range("B7").Value = "manchester@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:manchester@gmail.com", "", "", "manchester@gmail.com")

# --- Glasgow tribunal block (rows 8-12) ---
$ws.Range("A8").Value = "tribunalGlasgowAddress"
$ws.Range("B8").Value = "35 High Landing G3 6AD, Glasgow"

$ws.Range("A9").Value = "tribunalGlasgowTelephone"
$ws.Range("B9").Value = 3572531270

$ws.Range("A10").Value = "tribunalGlasgowFax"
$ws.Range("B10").Value = 2937126570

$ws.Range("A11").Value = "tribunalGlasgowDX"
$ws.Range("B11").Value = 1231123

$ws.Range("A12").Value = "tribunalGlasgowEmail"
$ws.Range("B12").Value = "glasgow@gmail.com"
$ws.Hyperlinks.Add($ws.Range("B12"), "mailto:glasgow@gmail.com", "", "", "glasgow@gmail.com")

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 29.5
$ws.Columns.Item(2).ColumnWidth = 29.83

# --- Selection moves on to the next empty row, as in the authored workbook ---
$ws.Range("B13").Select()
